# "Generate Report for Handback"
# The handback transform failed for the 0434d2f0 record because the
# handback file name did not match the expected handoff file name. Update
# the status + error detail for that record on both locale sheets, and
# widen the Error Detail column so the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcnError = "Handback file name: egjqshxa.x0l is different with handoff file name: 0434d2f0-b97f-42de-8fbf-4b4cd6d9412a.ab66390f717b508a20d27bcc668cf8b3cbbfb50c.zh-cn."
$dedeError = "Handback file name: egjqshxa.x0l is different with handoff file name: 0434d2f0-b97f-42de-8fbf-4b4cd6d9412a.ab66390f717b508a20d27bcc668cf8b3cbbfb50c.de-de."

# Row 3 on both locale sheets corresponds to the 0434d2f0-... record. Its
# Status flips from "Ready for handoff" to "Handback transform failed",
# and the Error Detail column records why the handback was rejected.
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("P3").Value = $zhcnError

$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("P3").Value = $dedeError

# The Overview sheet mirrors the same status text for the 0434d2f0-...
# record in its zh-cn/de-de columns.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# Widen the Error Detail column (P) so the new message is readable, same
# as the other wide text columns on these sheets. ColumnWidth is stored in
# "characters" and gets snapped to a pixel grid on write, so 39.1666...
# (39 + 1/6) is the value that round-trips to a stored width of exactly 40.
$zhcn.Range("P1").ColumnWidth = 39.16666666666667
$dede.Range("P1").ColumnWidth = 39.16666666666667
